$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A34").Value = "Nicholas Marzadro"
$ws.Range("B34").Value = "Matteo Zanlucchi | SBARX"
$ws.Range("C34").Value = "Matteo Diener | U.SGUARNA"
$ws.Range("D34").Value = "Federico Manica | IMONTAGNA"
$ws.Range("E34").Value = "Filippo Benetti | I Magnifici"
$ws.Range("F34").Value = "Alessandro Fanti | FC SALAGIARDINI"
